$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the DCA column (D) entirely; the "Green" values move into column C.
$ws.Columns.Item(4).Delete()

# Remove the two extra data rows (old rows 5 and 6); delete bottom-up so
# row indices of rows still to be removed don't shift.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Update header row.
$ws.Range("A1").Value = "Project details"
$ws.Range("B1").Value = "WLC (forecast)"
$ws.Range("C1").Value = "DCA"

# Update data rows with new project details (note embedded line breaks).
$ws.Range("A2").Value = "A11," + [char]10 + "£90m"
$ws.Range("B2").Value = 89
$ws.Range("C2").Value = "Green"

$ws.Range("A3").Value = "Columbia," + [char]10 + "£4,3bn"
$ws.Range("B3").Value = 4345
$ws.Range("C3").Value = "Green"

$ws.Range("A4").Value = "A13," + [char]10 + "£89,8bn"
$ws.Range("B4").Value = 89809
$ws.Range("C4").Value = "Green"
